$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Число спорт. сооруж." row under "Уровень жизни" table:
# D11 already holds "Жил. площ.на одного чел. - livarea ..." — just
# restyle it to match the regular bordered cell style (same as C11/B11),
# which drops the special "no right border" styling it had before.
$ws.Range("C11").Copy()
$ws.Range("D11").PasteSpecial(-4122)  # xlPasteFormats

# D12 gets the new value, with the same style as C12.
$ws.Range("D12").Value = "Число спорт. сооруж. - sportsvenue (шт.) (id8003001)"
$ws.Range("C12").Copy()
$ws.Range("D12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the active selection to match the saved workbook state.
$ws.Range("E15").Select()
